# Applies crypto price/volume updates from the GitHub Actions data refresh.
# Uses a scratch cell + Copy/PasteSpecial(values) round-trip so that numeric-
# looking strings (e.g. "1.00", "166.92") are written as literal text, just
# like the original inline-string cells, instead of being auto-coerced to
# numbers by the normal Range.Value assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")

function Set-TextValue([string]$addr, [string]$text) {
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

Set-TextValue "D2" '67.710.93'
Set-TextValue "E2" '  -1.79%  '
Set-TextValue "D3" '3.795.06'
Set-TextValue "E3" '  -1.22%  '
Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  +0.08%  '
Set-TextValue "D5" '596.11'
Set-TextValue "E5" '  -0.41%  '
Set-TextValue "D6" '166.92'
Set-TextValue "E6" '  -2.19%  '
Set-TextValue "D7" '3.793.10'
Set-TextValue "E7" '  -1.22%  '
Set-TextValue "E8" '  +0.04%  '
Set-TextValue "D9" '0.521'
Set-TextValue "E9" '  -0.60%  '
Set-TextValue "E10" '  -1.64%  '
Set-TextValue "D11" '6.36'
Set-TextValue "E11" '  -0.71%  '
Set-TextValue "E12" '  -1.32%  '
Set-TextValue "E13" '  -3.86%  '
Set-TextValue "D14" '36.00'
Set-TextValue "E14" '  -2.03%  '
Set-TextValue "D15" '4.432.95'
Set-TextValue "E15" '  -0.54%  '
Set-TextValue "D16" '3.796.11'
Set-TextValue "E16" '  -0.62%  '
Set-TextValue "D17" '18.59'
Set-TextValue "E17" '  +2.16%  '
Set-TextValue "D18" '67.680.88'
Set-TextValue "E18" '  -1.78%  '
Set-TextValue "D19" '7.08'
Set-TextValue "E19" '  +0.51%  '
Set-TextValue "E20" '  -0.23%  '
Set-TextValue "B21" 'BitcoinCash'
Set-TextValue "C21" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D21" '460.16'
Set-TextValue "E21" '  -2.05%  '
Set-TextValue "B22" 'Uniswap'
Set-TextValue "C22" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D22" '9.95'
Set-TextValue "E22" '  -9.52%  '
Set-TextValue "E23" '  -1.06%  '
Set-TextValue "E24" '  -2.54%  '
Set-TextValue "D25" '83.35'
Set-TextValue "E25" '  -1.42%  '
Set-TextValue "E26" '  +0.24%  '
Set-TextValue "E27" '  -4.82%  '
Set-TextValue "E28" '  -0.08%  '
Set-TextValue "D29" '9.98'
Set-TextValue "E29" '  -2.71%  '
Set-TextValue "D30" '3.942.52'
Set-TextValue "E30" '  -1.07%  '
Set-TextValue "E31" '  -0.65%  '
Set-TextValue "D32" '2.24'
Set-TextValue "E32" '  +1.89%  '
Set-TextValue "D33" '7.22'
Set-TextValue "E33" '  -2.73%  '
Set-TextValue "D34" '29.61'
Set-TextValue "E34" '  -2.71%  '
Set-TextValue "E35" '  -0.02%  '
Set-TextValue "D36" '9.05'
Set-TextValue "E36" '  -2.33%  '
Set-TextValue "D37" '0.0998'
Set-TextValue "E37" '  -1.83%  '
Set-TextValue "E38" '  -5.24%  '
Set-TextValue "E39" '  -1.25%  '
Set-TextValue "D40" '0.994'
Set-TextValue "E40" '  -1.08%  '
Set-TextValue "E41" '  -1.36%  '
Set-TextValue "D42" '1.00'
Set-TextValue "E42" '  +0.02%  '
Set-TextValue "D44" '48.15'
Set-TextValue "E44" '  +2.56%  '
Set-TextValue "D45" '43.84'
Set-TextValue "E45" '  -1.25%  '
Set-TextValue "D46" '0.296'
Set-TextValue "E46" '  -3.01%  '
Set-TextValue "D47" '150.94'
Set-TextValue "E47" '  +2.70%  '
Set-TextValue "D48" '8.30'
Set-TextValue "E48" '  -2.19%  '
Set-TextValue "D49" '27.15'
Set-TextValue "E49" '  +3.49%  '
Set-TextValue "D50" '390.87'
Set-TextValue "E50" '  -0.99%  '
Set-TextValue "E51" '  -6.81%  '

$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false

Write-Output "Applied crypto data update."
